$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "Sheet" to "Sheet1"
$ws.Name = "Sheet1"

# Style the header row (A1:D1): bold font, thin border all around, centered/top aligned
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Populate "Appreciated" column (D2:D32) with each researcher's common arXiv articles
$ws.Range("D2").Value = "arXiv:2310.02113, arXiv:2310.05269, arXiv:2310.11730, arXiv:2310.09665, arXiv:2310.13424, arXiv:2311.02100, arXiv:2310.01689, arXiv:2310.01676, arXiv:2310.06338, arXiv:2310.04975"
$ws.Range("D3").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"
$ws.Range("D4").Value = "arXiv:2310.02357, arXiv:2310.14261, arXiv:2310.01551, arXiv:2310.00526, arXiv:2310.05309, arXiv:2310.19647, arXiv:2311.05511, arXiv:2311.01927, arXiv:2310.04218, arXiv:2310.04425"
$ws.Range("D5").Value = "arXiv:2310.04022, arXiv:2310.10524, arXiv:2310.03557, arXiv:2310.03339, arXiv:2310.03254, arXiv:2310.03193, arXiv:2310.01283, arXiv:2310.00394, arXiv:2310.01046, arXiv:2310.00267"
$ws.Range("D6").Value = "arXiv:2310.02357, arXiv:2310.00603, arXiv:2310.07086, arXiv:2310.05592, arXiv:2310.09736, arXiv:2310.08977, arXiv:2310.15799, arXiv:2310.14261, arXiv:2311.04925, arXiv:2311.02802"
$ws.Range("D7").Value = "arXiv:2310.03598, arXiv:2310.03994, arXiv:2310.12254, arXiv:2310.11651, arXiv:2310.10653, arXiv:2310.18820, arXiv:2311.05498, arXiv:2311.05462, arXiv:2311.05037, arXiv:2311.02378"
$ws.Range("D8").Value = "arXiv:2310.00254, arXiv:2310.04975, arXiv:2310.01689, arXiv:2310.01676, arXiv:2311.02093, arXiv:2310.01594, arXiv:2309.16707, arXiv:2310.07471, arXiv:2310.08822, arXiv:2310.12381"
$ws.Range("D9").Value = "arXiv:2310.01557, arXiv:2310.03086, arXiv:2310.03243, arXiv:2310.02870, arXiv:2310.02422, arXiv:2310.00727, arXiv:2310.00633, arXiv:2310.00010, arXiv:2310.01063, arXiv:2309.16733"
$ws.Range("D10").Value = "arXiv:2311.00724, arXiv:2311.04517, arXiv:2311.04482, arXiv:2310.02113, arXiv:2310.01689, arXiv:2310.01676, arXiv:2310.06338, arXiv:2310.05269, arXiv:2310.04975, arXiv:2310.03841"
$ws.Range("D11").Value = "arXiv:2310.08373, arXiv:2310.14283, arXiv:2310.18861, arXiv:2310.09193, arXiv:2310.03618, arXiv:2310.03616, arXiv:2310.02800, arXiv:2310.01893, arXiv:2310.05643, arXiv:2310.12702"
$ws.Range("D12").Value = "arXiv:2310.00029, arXiv:2311.02389, arXiv:2310.01039, arXiv:2310.05170, arXiv:2310.15887, arXiv:2311.04569, arXiv:2311.04126, arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655"
$ws.Range("D13").Value = "arXiv:2310.01003, arXiv:2310.04870, arXiv:2310.03817, arXiv:2310.11678, arXiv:2310.13897, arXiv:2311.00208, arXiv:2311.02433, arXiv:2310.15388, arXiv:2310.10513, arXiv:2311.02460"
$ws.Range("D14").Value = "arXiv:2310.00382, arXiv:2311.04256, arXiv:2310.03401, arXiv:2310.01689, arXiv:2310.01676, arXiv:2310.00254, arXiv:2310.18664, arXiv:2310.18382, arXiv:2310.04975, arXiv:2311.04944"
$ws.Range("D15").Value = "arXiv:2311.00974, arXiv:2310.19013, arXiv:2310.00560, arXiv:2310.06141, arXiv:2310.11957, arXiv:2310.09665, arXiv:2310.16547, arXiv:2311.00271, arXiv:2311.02525, arXiv:2310.03618"
$ws.Range("D16").Value = "arXiv:2310.00782, arXiv:2309.16682, arXiv:2310.05293, arXiv:2310.15419, arXiv:2310.14977, arXiv:2310.17990, arXiv:2311.04333, arXiv:2311.02811, arXiv:2310.02420, arXiv:2310.08339"
$ws.Range("D17").Value = "arXiv:2311.04517, arXiv:2311.00724, arXiv:2311.04482, arXiv:2310.01689, arXiv:2310.01676, arXiv:2310.00254, arXiv:2310.04975, arXiv:2311.02093, arXiv:2310.08439, arXiv:2310.05701"
$ws.Range("D18").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"
$ws.Range("D19").Value = "arXiv:2310.18382, arXiv:2310.03401, arXiv:2310.01689, arXiv:2310.01676, arXiv:2310.00254, arXiv:2310.04975, arXiv:2310.18664, arXiv:2311.04944, arXiv:2311.02926, arXiv:2311.02093"
$ws.Range("D20").Value = "arXiv:2309.17315, arXiv:2310.00762, arXiv:2310.00290, arXiv:2310.08447, arXiv:2310.10316, arXiv:2310.18078, arXiv:2311.00049, arXiv:2310.19548, arXiv:2310.18565, arXiv:2311.03772"
$ws.Range("D21").Value = "arXiv:2310.02357, arXiv:2310.14261, arXiv:2310.07874, arXiv:2310.08039, arXiv:2310.04878, arXiv:2310.03919, arXiv:2310.11088, arXiv:2310.16452, arXiv:2310.14079, arXiv:2310.13006"
$ws.Range("D22").Value = "arXiv:2310.00278, arXiv:2310.06857, arXiv:2310.16106, arXiv:2310.16195, arXiv:2310.15705, arXiv:2310.14283, arXiv:2311.05582, arXiv:2310.03744, arXiv:2310.03743, arXiv:2310.03742"
$ws.Range("D23").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"
$ws.Range("D24").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"
$ws.Range("D25").Value = "arXiv:2310.02422, arXiv:2311.02840, arXiv:2310.06916, arXiv:2310.04837, arXiv:2310.03294, arXiv:2310.00627, arXiv:2309.16743, arXiv:2310.08401, arXiv:2310.08097, arXiv:2310.07471"
$ws.Range("D26").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"
$ws.Range("D27").Value = "arXiv:2310.01039, arXiv:2310.05170, arXiv:2310.15887, arXiv:2311.04569, arXiv:2311.04126, arXiv:2310.00029, arXiv:2311.02389, arXiv:2310.03673, arXiv:2310.03659, arXiv:2310.03620"
$ws.Range("D28").Value = "arXiv:2310.16214, arXiv:2311.03373, arXiv:2310.04172, arXiv:2311.05063, arXiv:2310.03618, arXiv:2310.03616, arXiv:2310.03568, arXiv:2310.03371, arXiv:2310.03294, arXiv:2310.03200"
$ws.Range("D29").Value = "arXiv:2310.03673, arXiv:2310.03659, arXiv:2310.03620, arXiv:2310.03618, arXiv:2310.03616, arXiv:2310.03533, arXiv:2310.03491, arXiv:2310.03318, arXiv:2310.03248, arXiv:2310.03202"
$ws.Range("D30").Value = "arXiv:2310.05020, arXiv:2310.02540, arXiv:2310.02129, arXiv:2310.01765, arXiv:2310.00749, arXiv:2310.07875, arXiv:2310.07736, arXiv:2310.04830, arXiv:2310.04598, arXiv:2310.04145"
$ws.Range("D31").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"
$ws.Range("D32").Value = "arXiv:2310.03736, arXiv:2310.03702, arXiv:2310.03655, arXiv:2310.03528, arXiv:2310.03501, arXiv:2310.03475, arXiv:2310.03441, arXiv:2310.03178, arXiv:2310.03159, arXiv:2310.03105"

Write-Host "edit complete"
